# "cleaned defensive actions data"
# Restructure the header rows (shift short labels into a new hidden row1 +
# a visible row2 that gains a "Player ID" column), rename a couple of
# columns, fill in previously-blank Tkl% zeros, hide the blank separator
# row and the summary row, and drop the old header merge cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 (becomes a hidden row of short column headers) ----
$row1 = @("Player","#","Nation","Pos","Age","90s","Tkl","TklW","Def 3rd","Mid 3rd","Att 3rd","Cha","Att","Tkl%","Lost","Blocks","Sh","Pass","Int","Tkl+Int","Clr","Err","Unnamed: 21_level_0")
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $row1[$i]
}

# ---- Row 2 (visible header; gains "Player ID" in column A) ----
$row2 = @("Player ID","Player","#","Nation","Pos","Age","90s","Tkl","TklW","Def 3rd","Mid 3rd","Att 3rd","Cha","Att","Tkl%","Lost","Blocks","Sh","Pass","Int","Tkl+Int","Clr","Err")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# ---- Fill previously-blank Tkl% cells (column O) with 0 ----
$zeroRows = @(4,5,7,9,15,16,19,20)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 15).Value = 0
}

# ---- Remove the old header merge cells ----
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# ---- Hide row 1 (now auxiliary), the blank row 3, and the summary row 21 ----
$ws.Rows.Item(1).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(21).Hidden = $true

Write-Host "defensive actions sheet cleaned"
